$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.883.03"
$ws.Range("E2").Value = "  +0.86%  "
$ws.Range("D3").Value = "2.638.41"
$ws.Range("E3").Value = "  +1.64%  "
$ws.Range("E4").Value = "  -0.13%  "
$ws.Range("D5").Value = "'578.32"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.69%  "
$ws.Range("D6").Value = "'144.10"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.90%  "
$ws.Range("E7").Value = "  +0.14%  "
$ws.Range("E8").Value = "  -0.32%  "
$ws.Range("E9").Value = "  +0.60%  "
$ws.Range("E10").Value = "  +2.13%  "
$ws.Range("E12").Value = "  +1.23%  "
$ws.Range("D13").Value = "3.106.72"
$ws.Range("E13").Value = "  +1.14%  "
$ws.Range("D14").Value = "'26.19"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +6.85%  "
$ws.Range("D15").Value = "60.845.24"
$ws.Range("E15").Value = "  +0.77%  "
$ws.Range("E16").Value = "  +1.82%  "
$ws.Range("D17").Value = "2.649.42"
$ws.Range("D18").Value = "'11.59"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.27%  "
$ws.Range("D19").Value = "'4.72"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.07%  "
$ws.Range("D20").Value = "'351.47"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.08%  "
$ws.Range("E21").Value = "  -0.54%  "
$ws.Range("D22").Value = "'0.999"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.06%  "
$ws.Range("D23").Value = "'0.526"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.13%  "
$ws.Range("D24").Value = "'63.92"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.39%  "
$ws.Range("E25").Value = "  +1.49%  "
$ws.Range("D26").Value = "'0.993"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.87%  "
$ws.Range("D27").Value = "'8.38"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +5.82%  "
$ws.Range("D28").Value = "'1.99"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +7.28%  "
$ws.Range("D29").Value = "0.0₃0808"
$ws.Range("E29").Value = "  +1.66%  "
$ws.Range("D30").Value = "'6.74"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +5.89%  "
$ws.Range("D31").Value = "'167.40"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.89%  "
$ws.Range("E32").Value = "  +0.03%  "
$ws.Range("D33").Value = "'19.94"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.27%  "
$ws.Range("E34").Value = "  +8.15%  "
$ws.Range("E35").Value = "  +8.93%  "
$ws.Range("E36").Value = "  +7.09%  "
$ws.Range("E37").Value = "  +4.29%  "
$ws.Range("D38").Value = "'341.28"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +9.79%  "
$ws.Range("E39").Value = "  +5.66%  "
$ws.Range("E40").Value = "  +7.24%  "
$ws.Range("D41").Value = "'38.22"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.74%  "
$ws.Range("D42").Value = "'138.36"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.06%  "
$ws.Range("D43").Value = "'5.30"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +4.18%  "
$ws.Range("D44").Value = "'21.08"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +4.60%  "
$ws.Range("E45").Value = "  +3.87%  "
$ws.Range("E46").Value = "  +3.66%  "
$ws.Range("E47").Value = "  +2.29%  "
$ws.Range("E48").Value = "  +3.64%  "
$ws.Range("D49").Value = "'0.0995"
$ws.Range("D49").Style = "Normal"
$ws.Range("E50").Value = "  +0.02%  "
$ws.Range("D51").Value = "2.086.22"
$ws.Range("E51").Value = "  +2.43%  "
